$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.802135944366455
$ws.Range("B1").Value = 5.021834850311279
$ws.Range("C1").Value = 3.618614912033081
$ws.Range("D1").Value = 2.103127956390381
$ws.Range("E1").Value = 1.852848768234253
